$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace every occurrence of the literal value "Gasoline" with "Fuel"
# across the used range (columns C and D hold the commodity names).
$used = $ws.UsedRange
$used.Replace("Gasoline", "Fuel", 1) | Out-Null

# Reflect the new active selection recorded in the edited workbook.
$ws.Range("L14").Select() | Out-Null
